# Auto-generated Excel COM-interop script to apply the Ridill_Profits market-data refresh
# (scheduled runner updating currentAveragePrice* / Leve* profit columns per sheet)
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 117.5
$ws.Range("I55").Value = 103.57143
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 103.57143
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = 110.42857
$ws.Range("N55").Value = -578
$ws.Range("H64").Value = 2918.476
$ws.Range("I64").Value = 2915.7896
$ws.Range("J64").Value = 2944
$ws.Range("K64").Value = 2915.7896
$ws.Range("L64").Value = 2944
$ws.Range("M64").Value = -2667.7896
$ws.Range("N64").Value = -3440
$ws.Range("H67").Value = 2918.476
$ws.Range("I67").Value = 2915.7896
$ws.Range("J67").Value = 2944
$ws.Range("K67").Value = 2915.7896
$ws.Range("L67").Value = 2944
$ws.Range("M67").Value = -2057.7896
$ws.Range("N67").Value = -4660
$ws.Range("H70").Value = 3519.6155
$ws.Range("I70").Value = 994.5
$ws.Range("J70").Value = 4641.8887
$ws.Range("K70").Value = 2983.5
$ws.Range("L70").Value = 13925.6661
$ws.Range("M70").Value = -2713.5
$ws.Range("N70").Value = -14465.6661
$ws.Range("H73").Value = 3519.6155
$ws.Range("I73").Value = 994.5
$ws.Range("J73").Value = 4641.8887
$ws.Range("K73").Value = 2983.5
$ws.Range("L73").Value = 13925.6661
$ws.Range("M73").Value = -2047.5
$ws.Range("N73").Value = -15797.6661
$ws.Range("H74").Value = 4046.5
$ws.Range("I74").Value = 3662.7778
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 3662.7778
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -2726.7778
$ws.Range("N74").Value = -9372
$ws.Range("H76").Value = 166671920
$ws.Range("I76").Value = 200005740
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 200005740
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -200005425
$ws.Range("N76").Value = -3430
$ws.Range("H77").Value = 4046.5
$ws.Range("I77").Value = 3662.7778
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 18313.889
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -13633.889
$ws.Range("N77").Value = -46860
$ws.Range("H79").Value = 166671920
$ws.Range("I79").Value = 200005740
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 200005740
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -200004648
$ws.Range("N79").Value = -4984
$ws.Range("H121").Value = 766.55554
$ws.Range("I121").Value = 488.57144
$ws.Range("J121").Value = 833.65515
$ws.Range("K121").Value = 1465.71432
$ws.Range("L121").Value = 2500.96545
$ws.Range("M121").Value = 281.28568
$ws.Range("N121").Value = -5994.96545

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2154
$ws.Range("I122").Value = 2154
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6462
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4012
$ws.Range("N122").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1687.2727
$ws.Range("I105").Value = 1676
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 1676
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = 71
$ws.Range("N105").Value = -5294
$ws.Range("H134").Value = 8956135
$ws.Range("I134").Value = 11531485
$ws.Range("J134").Value = 41464
$ws.Range("K134").Value = 34594455
$ws.Range("L134").Value = 124392
$ws.Range("M134").Value = -34591920
$ws.Range("N134").Value = -129462

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1941178
$ws.Range("I31").Value = 3334590
$ws.Range("K31").Value = 3334590
$ws.Range("M31").Value = -3334295
$ws.Range("H34").Value = 1941178
$ws.Range("I34").Value = 3334590
$ws.Range("K34").Value = 3334590
$ws.Range("M34").Value = -3334388
$ws.Range("H55").Value = 11500
$ws.Range("I55").Value = 7250
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 7250
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -6935
$ws.Range("N55").Value = -20630
$ws.Range("H62").Value = 2928.1538
$ws.Range("I62").Value = 2280
$ws.Range("K62").Value = 2280
$ws.Range("M62").Value = -1656
$ws.Range("H65").Value = 2928.1538
$ws.Range("I65").Value = 2280
$ws.Range("K65").Value = 11400
$ws.Range("M65").Value = -8280

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2101.9697
$ws.Range("I34").Value = 653.53845
$ws.Range("J34").Value = 3043.45
$ws.Range("K34").Value = 1960.61535
$ws.Range("L34").Value = 9130.349999999999
$ws.Range("M34").Value = -1876.61535
$ws.Range("N34").Value = -9298.349999999999
$ws.Range("H64").Value = 4151.727
$ws.Range("I64").Value = 636.25
$ws.Range("J64").Value = 4932.9443
$ws.Range("K64").Value = 1908.75
$ws.Range("L64").Value = 14798.8329
$ws.Range("M64").Value = -1638.75
$ws.Range("N64").Value = -15338.8329
$ws.Range("H67").Value = 4151.727
$ws.Range("I67").Value = 636.25
$ws.Range("J67").Value = 4932.9443
$ws.Range("K67").Value = 1908.75
$ws.Range("L67").Value = 14798.8329
$ws.Range("M67").Value = -972.75
$ws.Range("N67").Value = -16670.8329

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5871523
$ws.Range("I70").Value = 2020425
$ws.Range("J70").Value = 35717536
$ws.Range("K70").Value = 2020425
$ws.Range("L70").Value = 35717536
$ws.Range("M70").Value = -2020155
$ws.Range("N70").Value = -35718076
$ws.Range("H73").Value = 5871523
$ws.Range("I73").Value = 2020425
$ws.Range("J73").Value = 35717536
$ws.Range("K73").Value = 2020425
$ws.Range("L73").Value = 35717536
$ws.Range("M73").Value = -2019489
$ws.Range("N73").Value = -35719408
$ws.Range("H80").Value = 8150.909
$ws.Range("I80").Value = 4300
$ws.Range("J80").Value = 15852.728
$ws.Range("K80").Value = 4300
$ws.Range("L80").Value = 15852.728
$ws.Range("M80").Value = -3302
$ws.Range("N80").Value = -17848.728
$ws.Range("H83").Value = 8150.909
$ws.Range("I83").Value = 4300
$ws.Range("J83").Value = 15852.728
$ws.Range("K83").Value = 21500
$ws.Range("L83").Value = 79263.64
$ws.Range("M83").Value = -16508
$ws.Range("N83").Value = -89247.64

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 12386.5
$ws.Range("I45").Value = 500
$ws.Range("J45").Value = 16348.667
$ws.Range("K45").Value = 500
$ws.Range("L45").Value = 16348.667
$ws.Range("M45").Value = -93
$ws.Range("N45").Value = -17162.667
$ws.Range("H61").Value = 2357.926
$ws.Range("I61").Value = 1804.1111
$ws.Range("J61").Value = 3465.5557
$ws.Range("K61").Value = 1804.1111
$ws.Range("L61").Value = 3465.5557
$ws.Range("M61").Value = -1602.1111
$ws.Range("N61").Value = -3869.5557
$ws.Range("H97").Value = 19000
$ws.Range("J97").Value = 19000
$ws.Range("L97").Value = 19000
$ws.Range("N97").Value = -20982
$ws.Range("H113").Value = 2357.926
$ws.Range("I113").Value = 1804.1111
$ws.Range("J113").Value = 3465.5557
$ws.Range("K113").Value = 1804.1111
$ws.Range("L113").Value = 3465.5557
$ws.Range("M113").Value = 365.8888999999999
$ws.Range("N113").Value = -7805.5557

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 68194776
$ws.Range("I62").Value = 187527630
$ws.Range("J62").Value = 4571.4287
$ws.Range("K62").Value = 187527630
$ws.Range("L62").Value = 4571.4287
$ws.Range("M62").Value = -187527006
$ws.Range("N62").Value = -5819.4287
$ws.Range("H65").Value = 68194776
$ws.Range("I65").Value = 187527630
$ws.Range("J65").Value = 4571.4287
$ws.Range("K65").Value = 937638150
$ws.Range("L65").Value = 22857.1435
$ws.Range("M65").Value = -937635030
$ws.Range("N65").Value = -29097.1435
$ws.Range("H113").Value = 323.6087
$ws.Range("I113").Value = 242.78572
$ws.Range("J113").Value = 449.33334
$ws.Range("K113").Value = 728.35716
$ws.Range("L113").Value = 1348.00002
$ws.Range("M113").Value = 1441.64284
$ws.Range("N113").Value = -5688.000019999999

Write-Output "Applied all Ridill_Profits market-data updates"